# lsn 10 summary + lsn 10 grades + updated grades image for last two lsns
#
# Adds the new row of grades/attendance (row 26) for lesson 12, mirroring
# the layout already used by the other weekly rows (C/D/E/I/J/K/L/M hold the
# "attendance: 5" marker, N holds the lesson title), and updates the view so
# the newly entered row is visible/selected, the way it would be after a
# user had just typed it in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$attendance = "الحضور: 5"
$lesson12 = "الدرس 12 (أخلاق فتية حول الرسول 2)"

$row = 26
foreach ($col in @("C", "D", "E", "I", "J", "K", "L", "M")) {
    $ws.Range("$col$row").Value = $attendance
}
$ws.Range("N$row").Value = $lesson12

# Reflect the scrolled/selected state that results from having just added
# this row (frozen header pane now showing rows starting further down, and
# the active cell sitting just past the newly-added row).
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1

$ws.Range("N28").Select()
